$d = $word.ActiveDocument
$p1 = $d.Paragraphs(1)

# Replace the placeholder ID text in the first paragraph's run.
$p1.Range.Find.Execute("**ID__AFFARS_pgi_5301_topic_30__ID**", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "**ID__AFFARS_AFMC_PGI_5301_603_90__ID**", 2)

# Remove the trailing " " run that followed the placeholder text.
$p1.Range.Find.Execute("**ID__AFFARS_AFMC_PGI_5301_603_90__ID** ", $true, $false, $false, $false, $false, `
                        $true, 1, $false, "**ID__AFFARS_AFMC_PGI_5301_603_90__ID**", 2)

# Update the left indent from 120 twips (6pt) to 225 twips (11.25pt).
$p1.Format.LeftIndent = 11.25

# Add a paragraph border (no visible line, just the 5-twip spacing), matching
# the pBdr already used elsewhere in the document.
$p1.Borders.DistanceFromTop = 5
$p1.Borders.DistanceFromLeft = 5
$p1.Borders.DistanceFromBottom = 5
$p1.Borders.DistanceFromRight = 5
